# Journal de bord - "Updated with work done today"
# Fills in the 5 previously-empty rows (63-67) of the activity log with the
# day's work, and extends the trailing blank rows by one (150 -> 151), plus
# a couple of small incidental cell additions that Excel produced on save
# (E133 and D144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 64: "Réalisation" / endpoint de sélection d'activité ---
$ws.Range("A64").Value = "Réalisation"
$ws.Range("B64").Value = "Création du endpoint de sélection d'activité"
$ws.Range("C64").Value = 1.5
$ws.Range("D64").Value = 43538

# --- Row 65: "Réalisation" / endpoint de création + sélection du type d'activité ---
$ws.Range("A65").Value = "Réalisation"
$ws.Range("B65").Value = "Création du endpoint de création de type d'activité ainsi que de sélection du type d'activité"
$ws.Range("C65").Value = 1.5
$ws.Range("D65").Value = 43538

# --- Row 66: "Réalisation" / page de login admin ---
$ws.Range("A66").Value = "Réalisation"
$ws.Range("B66").Value = "Création de la page de login de l'interface web d'administration à l'aide des maquettes graphique réalisé précédement"
$ws.Range("C66").Value = 1
$ws.Range("D66").Value = 43538

# --- Row 67: "Réalisation" / déploiement page de login ---
$ws.Range("A67").Value = "Réalisation"
$ws.Range("B67").Value = "Déploiement de la page de login sur le serveur node, déplacement de tout l'api et adalptation du code pour prendre en compte les changements"
$ws.Range("C67").Value = 1.5
$ws.Range("D67").Value = 43538

# --- Row 63: "Réalisation" / continuation debug endpoint d'activité ---
$ws.Range("A63").Value = "Réalisation"
$ws.Range("B63").Value = "Continuation du debug de l'endpoint d'activité, des bugs étaient présent lors de l'insertion d'une nouvelle activité car la contrainte de clé étrangère du type d'activité n'était pas satisfaite"
$ws.Range("C63").Value = 1
$ws.Range("D63").Value = 43538

# Column A had no cells at all on rows 63-67 previously, so the new cells
# come in with the default style; apply the same wrapped-text style ("s=1")
# used by every other populated cell in column A/B/E.
$ws.Range("A63:A67").WrapText = $true

# Wrapped-text row heights (best-effort match of Excel's own auto-fit result
# for these cell contents).
$ws.Rows.Item(63).RowHeight = 105
$ws.Rows.Item(64).RowHeight = 30
$ws.Rows.Item(65).RowHeight = 60
$ws.Rows.Item(66).RowHeight = 60
$ws.Rows.Item(67).RowHeight = 90

# Keep the selection cursor where the author left off (one row further down).
$ws.Range("D64").Select()

# A couple of incidental style-only cells Excel also wrote further down the
# sheet (no value, just inherit the same formatting as their row neighbours).
$ws.Range("E132").Copy()
$ws.Range("E133").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D144").PasteSpecial(-4122)

# Extend the trailing empty rows by one (150 -> 151), matching the existing
# pattern used by the other blank rows at the bottom of the log.
$ws.Range("B150").Copy()
$ws.Range("B151").PasteSpecial(-4122)

$excel.CutCopyMode = $false
